$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 39 summary-statistics formulas.
# D39 ("Desvio Médio") now computes the MEDIAN instead of the average deviation.
$ws.Range("D39").Formula = "=MEDIAN(B1:B39)"

# E39 ("Desvio Padrão") becomes a plain placeholder value (no formula anymore).
$ws.Range("E39").Value = 0

# The standard-deviation and population-variance formulas shift one
# column to the right (from E/F into F/G).
$ws.Range("F39").Formula = "=STDEV.P(B1:B39)"
$ws.Range("G39").Formula = "=VAR.P(B1:B39)"

# H39 ("Variância Amostral") keeps its original formula/position.
$ws.Range("H39").Formula = "=VAR.S(B1:B39)"

# Update the active sheet view: scroll position and current selection.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("E39").Select()
